# "creando excel BOLSAS Total"
# GRUMA's row (row 27) moves from an explicit numeric "1" into the shared
# "ok" marker used by every other row in the F (estado/total) column. The
# dependent totals (F1 via =F40, and F40 via =SUM(F2:F39)) drop from 13 to
# 12 automatically once recalculated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F27: numeric 1 -> shared string "ok" (same text already used by F2:F26 / F28:F39)
$ws.Range("F27").Value = "ok"

# Reflect the cursor having moved one row down (F27 -> F28) after the edit,
# matching the saved selection state in the workbook.
$ws.Range("F28").Select()
